# Update the Germany Oberliga Hamburg match-odds table: several rows were
# re-ordered (their B:AC data — id through PL_AhUnder — moved to a
# different row while the leading row-index column A stayed put).
#
# We implement this by reading every source cell's value with .Value2
# (preserves the original numeric/string value) and writing it back to the
# destination cell, column by column, cell by cell (never through an
# array/range bulk assignment) so Excel's own cached representation of each
# value is reused intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-Rows($ws, $cols, $r1, $r2) {
    $vals1 = @()
    $vals2 = @()
    foreach ($c in $cols) {
        $vals1 += , $ws.Range("$c$r1").Value2
        $vals2 += , $ws.Range("$c$r2").Value2
    }
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r1").Value2 = $vals2[$i]
        $ws.Range("$($cols[$i])$r2").Value2 = $vals1[$i]
    }
}

# Simple pairwise swaps (A column / row index untouched)
Swap-Rows $ws $cols 38 39
Swap-Rows $ws $cols 50 51
Swap-Rows $ws $cols 57 58
Swap-Rows $ws $cols 59 61
Swap-Rows $ws $cols 263 264
Swap-Rows $ws $cols 275 276

# Three-way rotation among rows 288, 289, 290:
#   new288 = old289, new289 = old290, new290 = old288
$vals288 = @()
$vals289 = @()
$vals290 = @()
foreach ($c in $cols) {
    $vals288 += , $ws.Range("$c" + "288").Value2
    $vals289 += , $ws.Range("$c" + "289").Value2
    $vals290 += , $ws.Range("$c" + "290").Value2
}
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])288").Value2 = $vals289[$i]
    $ws.Range("$($cols[$i])289").Value2 = $vals290[$i]
    $ws.Range("$($cols[$i])290").Value2 = $vals288[$i]
}
